$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New procedure name row (added to the shared string table first so it
# occupies the earlier index, matching the target workbook layout)
$ws.Range("A6").Value = "sp_GetCurrentTimeForAgency"

# New "Add Procedures" header row (bold, mirrors the "Alter Procedures" header style)
$ws.Range("A5").Value = "Add Procedures"
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Font.Bold = $true

$ws.Range("A5").Select()
